$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Ceacam1"
$ws.Cells.Item(2,3).Value = "Sele"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 29.70764066666667
$ws.Cells.Item(2,8).Value = 89.122922
$ws.Cells.Item(2,9).Value = 0.8900806065804322
$ws.Cells.Item(2,10).Value = 0.8900806065804322
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 7.879565666666667
$ws.Cells.Item(2,14).Value = 23.638697
$ws.Cells.Item(2,15).Value = 0.9977172793687663
$ws.Cells.Item(2,16).Value = 0.9977172793687664
$ws.Cells.Item(2,17).Value = 234.0833054347371
$ws.Cells.Item(2,18).Value = 2106.749748912634
$ws.Cells.Item(2,19).Value = 0.88804880121633
$ws.Cells.Item(2,20).Value = 0.8880488012163301
# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Ceacam1"
$ws.Cells.Item(3,3).Value = "Sele"
$ws.Cells.Item(3,4).Value = "sCs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 29.70764066666667
$ws.Cells.Item(3,8).Value = 89.122922
$ws.Cells.Item(3,9).Value = 0.8900806065804322
$ws.Cells.Item(3,10).Value = 0.8900806065804322
$ws.Cells.Item(3,11).Value = 1
$ws.Cells.Item(3,12).Value = 0.3333333333333333
$ws.Cells.Item(3,13).Value = 0.018028
$ws.Cells.Item(3,14).Value = 0.054084
$ws.Cells.Item(3,15).Value = 0.002282720631233623
$ws.Cells.Item(3,16).Value = 0.002282720631233623
$ws.Cells.Item(3,17).Value = 0.5355693459386667
$ws.Cells.Item(3,18).Value = 4.820124113448
$ws.Cells.Item(3,19).Value = 0.002031805364102091
$ws.Cells.Item(3,20).Value = 0.002031805364102091
# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Ceacam1"
$ws.Cells.Item(4,3).Value = "Sele"
$ws.Cells.Item(4,4).Value = "ECs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 1.414758333333333
$ws.Cells.Item(4,8).Value = 4.244275
$ws.Cells.Item(4,9).Value = 0.04238804991710397
$ws.Cells.Item(4,10).Value = 0.04238804991710397
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 7.879565666666667
$ws.Cells.Item(4,14).Value = 23.638697
$ws.Cells.Item(4,15).Value = 0.9977172793687663
$ws.Cells.Item(4,16).Value = 0.9977172793687664
$ws.Cells.Item(4,17).Value = 11.14768118996389
$ws.Cells.Item(4,18).Value = 100.329130709675
$ws.Cells.Item(4,19).Value = 0.04229128984104043
$ws.Cells.Item(4,20).Value = 0.04229128984104044
# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Ceacam1"
$ws.Cells.Item(5,3).Value = "Sele"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 1.414758333333333
$ws.Cells.Item(5,8).Value = 4.244275
$ws.Cells.Item(5,9).Value = 0.04238804991710397
$ws.Cells.Item(5,10).Value = 0.04238804991710397
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.018028
$ws.Cells.Item(5,14).Value = 0.054084
$ws.Cells.Item(5,15).Value = 0.002282720631233623
$ws.Cells.Item(5,16).Value = 0.002282720631233623
$ws.Cells.Item(5,17).Value = 0.02550526323333333
$ws.Cells.Item(5,18).Value = 0.2295473691
$ws.Cells.Item(5,19).Value = 0.00009676007606353392
$ws.Cells.Item(5,20).Value = 0.00009676007606353392
# Row 6
$ws.Cells.Item(6,1).Value = "sCs"
$ws.Cells.Item(6,2).Value = "Ceacam1"
$ws.Cells.Item(6,3).Value = "Sele"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 2.253949666666667
$ws.Cells.Item(6,8).Value = 6.761849
$ws.Cells.Item(6,9).Value = 0.06753134350246381
$ws.Cells.Item(6,10).Value = 0.0675313435024638
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 7.879565666666667
$ws.Cells.Item(6,14).Value = 23.638697
$ws.Cells.Item(6,15).Value = 0.9977172793687663
$ws.Cells.Item(6,16).Value = 0.9977172793687664
$ws.Cells.Item(6,17).Value = 17.76014440786145
$ws.Cells.Item(6,18).Value = 159.841299670753
$ws.Cells.Item(6,19).Value = 0.06737718831139582
$ws.Cells.Item(6,20).Value = 0.0673771883113958
# Row 7
$ws.Cells.Item(7,1).Value = "sCs"
$ws.Cells.Item(7,2).Value = "Ceacam1"
$ws.Cells.Item(7,3).Value = "Sele"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 2.253949666666667
$ws.Cells.Item(7,8).Value = 6.761849
$ws.Cells.Item(7,9).Value = 0.06753134350246381
$ws.Cells.Item(7,10).Value = 0.0675313435024638
$ws.Cells.Item(7,11).Value = 1
$ws.Cells.Item(7,12).Value = 0.3333333333333333
$ws.Cells.Item(7,13).Value = 0.018028
$ws.Cells.Item(7,14).Value = 0.054084
$ws.Cells.Item(7,15).Value = 0.002282720631233623
$ws.Cells.Item(7,16).Value = 0.002282720631233623
$ws.Cells.Item(7,17).Value = 0.04063420459066667
$ws.Cells.Item(7,18).Value = 0.365707841316
$ws.Cells.Item(7,19).Value = 0.0001541551910679989
$ws.Cells.Item(7,20).Value = 0.0001541551910679988
